$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C2:C13 column held percentage-formatted values (style index referencing
# numFmt "0.0%" + the red "Aptos Narrow" font + border, via the "Porcentaje"
# cell style). The edit converts these to plain numbers using the same
# formatting already used by the neighbouring D:I columns (red font + border,
# default/general number format) and rewrites the values as whole numbers
# (old percentage value * 100).

# Copy the formatting already used in column D (style shared by D2:I13) onto
# C2:C13 so the same underlying cell-format record is reused instead of a
# brand new one being minted.
$ws.Range("D2:D13").Copy() | Out-Null
$ws.Range("C2:C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newValues = @(101, 51, 67, 47, 87, 82, 73, 68, 86, 71, 66, 42)
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

# The percentage cell style ("Porcentaje") is no longer used by any cell now
# that C2:C13 no longer reference it - remove the now-unused named style.
foreach ($style in $wb.Styles) {
    if ($style.Name -eq "Porcentaje") {
        $style.Delete()
    }
}

# Selection moved from A14 to C14.
$ws.Range("C14").Select()
